$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last three data rows (previously sending-cluster "MuSCs" -> now removed entirely,
# since only FAPs/MuSCs remain as sending clusters after the TPM update)
$ws.Rows("8:10").Delete()

# Row 2
$ws.Range("A2").Value2 = "FAPs"
$ws.Range("B2").Value2 = "Ifnb1"
$ws.Range("C2").Value2 = "Ifnar1"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 1
$ws.Range("F2").Value2 = 0.3333333333333333
$ws.Range("G2").Value2 = 0.001763
$ws.Range("H2").Value2 = 0.005289
$ws.Range("I2").Value2 = 0.6312209094163982
$ws.Range("J2").Value2 = 0.6312209094163981
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 12.944643
$ws.Range("N2").Value2 = 38.833929
$ws.Range("O2").Value2 = 0.3748900893017936
$ws.Range("P2").Value2 = 0.3748900893017936
$ws.Range("Q2").Value2 = 0.022821405609
$ws.Range("R2").Value2 = 0.205392650481
$ws.Range("S2").Value2 = 0.2366384631002729
$ws.Range("T2").Value2 = 0.2366384631002729

# Row 3
$ws.Range("A3").Value2 = "FAPs"
$ws.Range("B3").Value2 = "Ifnb1"
$ws.Range("C3").Value2 = "Ifnar1"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 1
$ws.Range("F3").Value2 = 0.3333333333333333
$ws.Range("G3").Value2 = 0.001763
$ws.Range("H3").Value2 = 0.005289
$ws.Range("I3").Value2 = 0.6312209094163982
$ws.Range("J3").Value2 = 0.6312209094163981
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 15.033452
$ws.Range("N3").Value2 = 45.100356
$ws.Range("O3").Value2 = 0.4353841324781401
$ws.Range("P3").Value2 = 0.4353841324781401
$ws.Range("Q3").Value2 = 0.026503975876
$ws.Range("R3").Value2 = 0.238535782884
$ws.Range("S3").Value2 = 0.2748235680483212
$ws.Range("T3").Value2 = 0.2748235680483211

# Row 4
$ws.Range("A4").Value2 = "FAPs"
$ws.Range("B4").Value2 = "Ifnb1"
$ws.Range("C4").Value2 = "Ifnar1"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("E4").Value2 = 1
$ws.Range("F4").Value2 = 0.3333333333333333
$ws.Range("G4").Value2 = 0.001763
$ws.Range("H4").Value2 = 0.005289
$ws.Range("I4").Value2 = 0.6312209094163982
$ws.Range("J4").Value2 = 0.6312209094163981
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 6.551073333333334
$ws.Range("N4").Value2 = 19.65322
$ws.Range("O4").Value2 = 0.1897257782200662
$ws.Range("P4").Value2 = 0.1897257782200662
$ws.Range("Q4").Value2 = 0.01154954228666667
$ws.Range("R4").Value2 = 0.10394588058
$ws.Range("S4").Value2 = 0.1197588782678041
$ws.Range("T4").Value2 = 0.119758878267804

# Row 5
$ws.Range("A5").Value2 = "MuSCs"
$ws.Range("B5").Value2 = "Ifnb1"
$ws.Range("C5").Value2 = "Ifnar1"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 1
$ws.Range("F5").Value2 = 0.3333333333333333
$ws.Range("G5").Value2 = 0.00103
$ws.Range("H5").Value2 = 0.00309
$ws.Range("I5").Value2 = 0.3687790905836018
$ws.Range("J5").Value2 = 0.3687790905836018
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 12.944643
$ws.Range("N5").Value2 = 38.833929
$ws.Range("O5").Value2 = 0.3748900893017936
$ws.Range("P5").Value2 = 0.3748900893017936
$ws.Range("Q5").Value2 = 0.01333298229
$ws.Range("R5").Value2 = 0.11999684061
$ws.Range("S5").Value2 = 0.1382516262015207
$ws.Range("T5").Value2 = 0.1382516262015207

# Row 6
$ws.Range("A6").Value2 = "MuSCs"
$ws.Range("B6").Value2 = "Ifnb1"
$ws.Range("C6").Value2 = "Ifnar1"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = 0.3333333333333333
$ws.Range("G6").Value2 = 0.00103
$ws.Range("H6").Value2 = 0.00309
$ws.Range("I6").Value2 = 0.3687790905836018
$ws.Range("J6").Value2 = 0.3687790905836018
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 15.033452
$ws.Range("N6").Value2 = 45.100356
$ws.Range("O6").Value2 = 0.4353841324781401
$ws.Range("P6").Value2 = 0.4353841324781401
$ws.Range("Q6").Value2 = 0.01548445556
$ws.Range("R6").Value2 = 0.13936010004
$ws.Range("S6").Value2 = 0.1605605644298189
$ws.Range("T6").Value2 = 0.1605605644298189

# Row 7
$ws.Range("A7").Value2 = "MuSCs"
$ws.Range("B7").Value2 = "Ifnb1"
$ws.Range("C7").Value2 = "Ifnar1"
$ws.Range("D7").Value2 = "MuSCs"
$ws.Range("E7").Value2 = 1
$ws.Range("F7").Value2 = 0.3333333333333333
$ws.Range("G7").Value2 = 0.00103
$ws.Range("H7").Value2 = 0.00309
$ws.Range("I7").Value2 = 0.3687790905836018
$ws.Range("J7").Value2 = 0.3687790905836018
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 6.551073333333334
$ws.Range("N7").Value2 = 19.65322
$ws.Range("O7").Value2 = 0.1897257782200662
$ws.Range("P7").Value2 = 0.1897257782200662
$ws.Range("Q7").Value2 = 0.006747605533333333
$ws.Range("R7").Value2 = 0.0607284498
$ws.Range("S7").Value2 = 0.06996689995226214
$ws.Range("T7").Value2 = 0.06996689995226214
